# Added horizontal flip to train transforms: new resnext50_32x4d_00 folds
# results (fold0..fold4 + folds summary row), plus filling in the
# previously-pending metrics for the resnext50_32x4d_00 runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style already used by the surrounding "resnext50_32x4d_00_fold*" rows -
# reuse it for every brand new row so the new rows look consistent with
# the rest of the table.
$rowStyle = $ws.Range("A18:L18").Style

# --- Row 15 (resnext50_32x4d_11) : fill in previously empty metrics ---
$ws.Cells.Item(15, 7).Value = 30        # G15 min at
$ws.Cells.Item(15, 9).Value = 1.58      # I15 min val loss
$ws.Cells.Item(15, 10).Value = "-"      # J15 public_score

# --- Row 16 (resnext50_32x4d_12) : fill in metrics + comment ---
$ws.Cells.Item(16, 9).Value = 1.66      # I16 min val loss
$ws.Cells.Item(16, 10).Value = "-"      # J16 public_score
$ws.Cells.Item(16, 13).Value = "Хуже, чем 10й запуск с Radam lr=0.001"   # M16 comments

# --- Row 17 (resnext50_32x4d_00_fold0) : fill in metrics ---
$ws.Cells.Item(17, 7).Value = 44        # G17 min at
$ws.Cells.Item(17, 9).Value = 1.48      # I17 min val loss
$ws.Cells.Item(17, 10).Value = 9.43     # J17 public_score

# --- Row 18 (resnext50_32x4d_00_fold1) : fill in metrics ---
$ws.Cells.Item(18, 7).Value = 32        # G18 min at
$ws.Cells.Item(18, 9).Value = 1.51      # I18 min val loss
$ws.Cells.Item(18, 10).Value = 9.19     # J18 public_score

# --- Row 19 (new) resnext50_32x4d_00_fold2 ---
$ws.Cells.Item(19, 1).Value = "resnext50_32x4d_00_fold2"
$ws.Cells.Item(19, 2).Value = "Adam"
$ws.Cells.Item(19, 3).Value = 0.001
$ws.Cells.Item(19, 4).Value = "ReduceOnPlateau(patience=5, factor=0.5)"
$ws.Cells.Item(19, 5).Value = 512
$ws.Cells.Item(19, 6).Value = 80
$ws.Cells.Item(19, 7).Value = 36
$ws.Cells.Item(19, 8).Value = "11:00"
$ws.Cells.Item(19, 9).Value = 1.48
$ws.Cells.Item(19, 11).Value = 1
$ws.Cells.Item(19, 12).Value = 0.485
$ws.Range("A19:L19").Style = $rowStyle

# --- Row 20 (new) resnext50_32x4d_00_fold3 ---
$ws.Cells.Item(20, 1).Value = "resnext50_32x4d_00_fold3"
$ws.Cells.Item(20, 2).Value = "Adam"
$ws.Cells.Item(20, 3).Value = 0.001
$ws.Cells.Item(20, 4).Value = "ReduceOnPlateau(patience=5, factor=0.5)"
$ws.Cells.Item(20, 5).Value = 512
$ws.Cells.Item(20, 6).Value = 80
$ws.Cells.Item(20, 7).Value = 35
$ws.Cells.Item(20, 8).Value = "11:00"
$ws.Cells.Item(20, 9).Value = 1.49
$ws.Cells.Item(20, 11).Value = 1
$ws.Cells.Item(20, 12).Value = 0.485
$ws.Range("A20:L20").Style = $rowStyle

# --- Row 21 (new) resnext50_32x4d_00_fold4 ---
$ws.Cells.Item(21, 1).Value = "resnext50_32x4d_00_fold4"
$ws.Cells.Item(21, 2).Value = "Adam"
$ws.Cells.Item(21, 3).Value = 0.001
$ws.Cells.Item(21, 4).Value = "ReduceOnPlateau(patience=5, factor=0.5)"
$ws.Cells.Item(21, 5).Value = 512
$ws.Cells.Item(21, 6).Value = 80
$ws.Cells.Item(21, 8).Value = "11:00"
$ws.Cells.Item(21, 11).Value = 1
$ws.Cells.Item(21, 12).Value = 0.485
$ws.Range("A21:L21").Style = $rowStyle

# --- Row 22 (new) resnext50_32x4d_00_folds (combined summary) ---
$ws.Cells.Item(22, 1).Value = "resnext50_32x4d_00_folds"
$ws.Cells.Item(22, 10).Value = 8.33
$ws.Cells.Item(22, 11).Value = 1
$ws.Cells.Item(22, 12).Value = 0.485
$ws.Range("A22").Style = $rowStyle
$ws.Range("L22").Style = $rowStyle

# --- keep the active-cell selection in sync with the new last row ---
$ws.Range("L22").Select() | Out-Null
